$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 of the test-case sheet (DPLKAKT008-022) is being refreshed with a
# new userid (32382 -> 31160), so the PREPARATION notes text and the
# USERID cell both need to reflect the new id.
$ws.Range("F2").Value = "Username : 31160;`nPassword : bni1234;`nTgl. Market : 22/01/2023;`nFile Excel : 22012023HargaPasarFixedIncome.xlsx"
$ws.Range("G2").Value = 31160

# FILE_EXCEL stays the same file name.
$ws.Range("O2").Value = "22012023HargaPasarFixedIncome.xlsx"

# Scroll the sheet view so column E is the left-most visible column
# (previously it was scrolled to column G).
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
